# Update data: 2025-10-29 18:50
# Adds a new "distance from Dma50" worksheet with Sma50 distance data for indices.

$wb = $excel.ActiveWorkbook

# Add the new sheet as the last tab in the workbook (after the current last sheet).
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "distance from Dma50"

# Header row
$ws.Range("A1").Value = "Icon"
$ws.Range("B1").Value = "Stock"
$ws.Range("C1").Value = "Distance From Sma50"

# Reuse the existing bold/centered/bordered header style (same one used by
# the other sheets) instead of constructing a brand-new style.
$styleSource = $wb.Worksheets.Item("1 Month Performance").Range("A1:C1")
$styleSource.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)

$ICON = "📈"
$stocks = @(
  @("NIFTYPSUBANK", 10.2033),
  @("NIFTYMETAL", 8.6247),
  @("NIFTYOILANDGAS", 6.396),
  @("NIFTYCOMMODITIES", 5.7207),
  @("CNXINFRA", 5.6012),
  @("CNXREALTY", 5.4493),
  @("NIFTYPVTBANK", 5.0059),
  @("BANKNIFTY", 4.9192),
  @("NIFTYFINSERVICE", 3.9783),
  @("NIFTYMIDCAP50", 3.9228),
  @("NIFTY", 3.7191),
  @("CNXENERGY", 3.706),
  @("CNXMIDCAP", 3.6313),
  @("NIFTY200", 3.5915),
  @("NIFTY100", 3.5759),
  @("NIFTY500", 3.3379),
  @("CNXSMALLCAP", 2.9205),
  @("NIFTY50VALUE20", 2.893),
  @("NIFTYCPSE", 2.837),
  @("CNXNIFTYJUNIOR", 2.8291),
  @("NIFTYHEALTHCARE", 2.162),
  @("CNXIT", 2.0641),
  @("NIFTYCONSUMPTION", 2.0245),
  @("CNXPHARMA", 1.573),
  @("NIFTYAUTO", 1.5538),
  @("NIFTYGROWSECT15", 1.5325),
  @("NIFTYFMCG", 1.3194),
  @("NIFTYCONSURDURBL", 0.4031),
  @("NIFTYMEDIA", -1.9217)
)

$row = 2
foreach ($item in $stocks) {
  $ws.Cells.Item($row, 1).Value = $ICON
  $ws.Cells.Item($row, 2).Value = $item[0]
  $ws.Cells.Item($row, 3).Value = $item[1]
  $row = $row + 1
}

Write-Output "Added sheet distance from Dma50 with $($stocks.Length) data rows."
